$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 11:22"

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# Update case-count statistics for several countries (rows unaffected by the
# "Isla de Man" reorder below)
Set-RowValues 13  @(24983, 1580, 5164, 17296, 1285, 283, 2523)
Set-RowValues 19  @(13011, 69, 5240, 7476, 266, 22, 295)
Set-RowValues 37  @(4228, 109, 1608, 2553, 72, 2, 67)
Set-RowValues 65  @(1124, 33, 128, 953, 34, 3, 43)
Set-RowValues 68  @(974, 13, 293, 677, 14, 0, 4)
Set-RowValues 75  @(764, 37, 55, 702, 21, 0, 7)
Set-RowValues 93  @(409, 9, 165, 222, 7, 0, 22)
Set-RowValues 102 @(299, 0, 16, 281, 4, 1, 2)

# "Isla de Man" moves up in the list (now right after "Kenia", before
# "Venezuela"), with fresh data. Venezuela and Guinea each shift down one
# row, keeping their own data intact; Martinica (row 121) is unaffected.
$ws.Cells.Item(118, 1).Value = "Isla de Man"
Set-RowValues 118 @(171, 13, 88, 82, 10, 0, 1)

$ws.Cells.Item(119, 1).Value = "Venezuela"
Set-RowValues 119 @(167, 0, 65, 93, 6, 0, 9)

$ws.Cells.Item(120, 1).Value = "Guinea"
Set-RowValues 120 @(164, 0, 5, 159, 0, 0, 0)
